$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3, shifting the existing row 3
# (ttLF3hdB / Yokohama F. Marinos vs Urawa Reds) down to row 4.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the "Haras El Hodood vs Smouha" match.
$ws.Range("A3").Value = "nJe0zKQr"
$ws.Range("B3").Value = "30/10/2024"
$ws.Range("C3").Value = "11:00"
$ws.Range("D3").Value = "EGYPT - PREMIER LEAGUE"
$ws.Range("E3").Value = "Haras El Hodood"
$ws.Range("F3").Value = "Smouha"

$row3Vals = @(3.8,2.75,2.2,4.4,1.85,2.9,1.13,5.1,1.55,2.3,2.6,1.44,1.6,2.27,2.1,1.65,8,18.5,13.5,60,45,60,5.1,5.5,17.5,120,900,5.5,9.25,9.5,22,23,45,5.4,23,32,150,200,500,2.25,7.6,90,3.9,12.5,24,55,110,400,51,51)
for ($i = 0; $i -lt $row3Vals.Length; $i++) {
    $ws.Cells.Item(3, 7 + $i).Value = $row3Vals[$i]
}

# Fix the two odds values that changed slightly on the shifted row (row 4).
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2.2
